# Auto-generated edit script: updates crypto price/volume data
# per commit 'Updated cryptos list on Wed Oct 25 08:50:23 UTC 2023 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '33.971.67'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.774.38'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -2.19%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('E5').Value = '  -1.32%  '
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('E9').Value = '  -0.56%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0654'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.028.21'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.22%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.774.41'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.84'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +6.08%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '33.989.07'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.621'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.10%  '
$ws.Range('E17').Value = '  -1.72%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '68.58'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '251.84'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.25%  '
$ws.Range('E20').Value = '  -1.51%  '
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('E22').Value = '  -1.96%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.19'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.40%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.13'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.10%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '156.44'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '16.35'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.45%  '
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.75'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.10%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0508'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('E32').Value = '  -0.99%  '
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.83'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.92%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.445.00'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -6.11%  '
$ws.Range('E36').Value = '  -2.10%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.623'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('E40').Value = '  -2.17%  '
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('E42').Value = '  -2.59%  '
$ws.Range('E43').Value = '  -5.06%  '
$ws.Range('E44').Value = '  -2.50%  '
$ws.Range('E45').Value = '  -2.09%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.928.94'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.05%  '
$ws.Range('E47').Value = '  +0.68%  '
$ws.Range('E48').Value = '  +0.24%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '11.81'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.50%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '97.35'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.05%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '49.22'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -6.43%  '
